# Update countries & provincias Spain
#
# Refresh the "Pais" COVID snapshot: new timestamp, refreshed case
# counts for several countries, and a handful of countries that swapped
# rank/row order as their totals changed (their row keeps its position,
# but the country name + stats shown there change accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "last updated" timestamp footer
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 14 de Septiembre de 2020 a las 13:49"

# Row 61: Suiza - refreshed stats
$ws.Cells.Item(61, 2).Value2 = 47436
$ws.Cells.Item(61, 3).Value2 = 257
$ws.Cells.Item(61, 4).Value2 = 38900
$ws.Cells.Item(61, 5).Value2 = 6513
$ws.Cells.Item(61, 6).Value2 = 0
$ws.Cells.Item(61, 7).Value2 = 2
$ws.Cells.Item(61, 8).Value2 = 2023

# Rows 72/73: Irlanda and Estado de Palestina swap rank order
$ws.Cells.Item(72, 1).Value2 = "Estado de Palestina"
$ws.Cells.Item(72, 2).Value2 = 31362
$ws.Cells.Item(72, 3).Value2 = 788
$ws.Cells.Item(72, 4).Value2 = 21406
$ws.Cells.Item(72, 5).Value2 = 9730
$ws.Cells.Item(72, 6).Value2 = 0
$ws.Cells.Item(72, 7).Value2 = 5
$ws.Cells.Item(72, 8).Value2 = 226

$ws.Cells.Item(73, 1).Value2 = "Irlanda"
$ws.Cells.Item(73, 2).Value2 = 30985
$ws.Cells.Item(73, 3).Value2 = 0
$ws.Cells.Item(73, 4).Value2 = 23364
$ws.Cells.Item(73, 5).Value2 = 5837
$ws.Cells.Item(73, 6).Value2 = 0
$ws.Cells.Item(73, 7).Value2 = 0
$ws.Cells.Item(73, 8).Value2 = 1784

# Rows 85-87: refreshed stats (no rank changes)
$ws.Cells.Item(85, 2).Value2 = 15827
$ws.Cells.Item(85, 3).Value2 = 36
$ws.Cells.Item(85, 4).Value2 = 13315
$ws.Cells.Item(85, 5).Value2 = 1860
$ws.Cells.Item(85, 6).Value2 = 0
$ws.Cells.Item(85, 7).Value2 = 4
$ws.Cells.Item(85, 8).Value2 = 652

$ws.Cells.Item(86, 2).Value2 = 15769
$ws.Cells.Item(86, 3).Value2 = 12
$ws.Cells.Item(86, 4).Value2 = 14411
$ws.Cells.Item(86, 5).Value2 = 1145
$ws.Cells.Item(86, 6).Value2 = 0
$ws.Cells.Item(86, 7).Value2 = 2
$ws.Cells.Item(86, 8).Value2 = 213

$ws.Cells.Item(87, 2).Value2 = 14306
$ws.Cells.Item(87, 3).Value2 = 26
$ws.Cells.Item(87, 4).Value2 = 10563
$ws.Cells.Item(87, 5).Value2 = 3446
$ws.Cells.Item(87, 6).Value2 = 0
$ws.Cells.Item(87, 7).Value2 = 0
$ws.Cells.Item(87, 8).Value2 = 297

# Rows 116-120: Uganda jumps ahead of Hong Kong, Congo, Nicaragua, Cabo Verde
$ws.Cells.Item(116, 1).Value2 = "Uganda"
$ws.Cells.Item(116, 2).Value2 = 4978
$ws.Cells.Item(116, 3).Value2 = 179
$ws.Cells.Item(116, 4).Value2 = 2317
$ws.Cells.Item(116, 5).Value2 = 2605
$ws.Cells.Item(116, 6).Value2 = 0
$ws.Cells.Item(116, 7).Value2 = 1
$ws.Cells.Item(116, 8).Value2 = 56

$ws.Cells.Item(117, 1).Value2 = "Hong Kong"
$ws.Cells.Item(117, 2).Value2 = 4972
$ws.Cells.Item(117, 3).Value2 = 14
$ws.Cells.Item(117, 4).Value2 = 4635
$ws.Cells.Item(117, 5).Value2 = 236
$ws.Cells.Item(117, 6).Value2 = 0
$ws.Cells.Item(117, 7).Value2 = 1
$ws.Cells.Item(117, 8).Value2 = 101

$ws.Cells.Item(118, 1).Value2 = "Congo"
$ws.Cells.Item(118, 2).Value2 = 4928
$ws.Cells.Item(118, 3).Value2 = 0
$ws.Cells.Item(118, 4).Value2 = 3887
$ws.Cells.Item(118, 5).Value2 = 953
$ws.Cells.Item(118, 6).Value2 = 0
$ws.Cells.Item(118, 7).Value2 = 0
$ws.Cells.Item(118, 8).Value2 = 88

$ws.Cells.Item(119, 1).Value2 = "Nicaragua"
$ws.Cells.Item(119, 2).Value2 = 4818
$ws.Cells.Item(119, 3).Value2 = 0
$ws.Cells.Item(119, 4).Value2 = 2913
$ws.Cells.Item(119, 5).Value2 = 1761
$ws.Cells.Item(119, 6).Value2 = 0
$ws.Cells.Item(119, 7).Value2 = 0
$ws.Cells.Item(119, 8).Value2 = 144

$ws.Cells.Item(120, 1).Value2 = "Cabo Verde"
$ws.Cells.Item(120, 2).Value2 = 4813
$ws.Cells.Item(120, 3).Value2 = 0
$ws.Cells.Item(120, 4).Value2 = 4119
$ws.Cells.Item(120, 5).Value2 = 650
$ws.Cells.Item(120, 6).Value2 = 0
$ws.Cells.Item(120, 7).Value2 = 0
$ws.Cells.Item(120, 8).Value2 = 44

# Rows 145/146: Malta and Georgia swap rank order
$ws.Cells.Item(145, 1).Value2 = "Malta"
$ws.Cells.Item(145, 2).Value2 = 2405
$ws.Cells.Item(145, 3).Value2 = 53
$ws.Cells.Item(145, 4).Value2 = 1890
$ws.Cells.Item(145, 5).Value2 = 499
$ws.Cells.Item(145, 6).Value2 = 0
$ws.Cells.Item(145, 7).Value2 = 1
$ws.Cells.Item(145, 8).Value2 = 16

$ws.Cells.Item(146, 1).Value2 = "Georgia"
$ws.Cells.Item(146, 2).Value2 = 2392
$ws.Cells.Item(146, 3).Value2 = 165
$ws.Cells.Item(146, 4).Value2 = 1369
$ws.Cells.Item(146, 5).Value2 = 1004
$ws.Cells.Item(146, 6).Value2 = 0
$ws.Cells.Item(146, 7).Value2 = 0
$ws.Cells.Item(146, 8).Value2 = 19

# Row 194: refreshed stats (no rank change)
$ws.Cells.Item(194, 2).Value2 = 140
$ws.Cells.Item(194, 3).Value2 = 1
$ws.Cells.Item(194, 4).Value2 = 136
$ws.Cells.Item(194, 5).Value2 = 4
$ws.Cells.Item(194, 6).Value2 = 0
$ws.Cells.Item(194, 7).Value2 = 0
$ws.Cells.Item(194, 8).Value2 = 0
